$wb = $excel.ActiveWorkbook

# --- Boolean sheet: split the BVTQaZ.csv and VTQaZ.csv rows into 6 mode-specific
#     rows each (LDVs, HDVs, aircraft, rail, ships, motorbikes). ---
$wsBool = $wb.Worksheets.Item("Boolean")

# BVTQaZ.csv (row 17) first.
$wsBool.Rows("18:22").Insert()
$wsBool.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$wsBool.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$wsBool.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$wsBool.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$wsBool.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$wsBool.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# VTQaZ.csv was row 21, now shifted down 5 rows to row 26.
$wsBool.Rows("27:31").Insert()
$wsBool.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$wsBool.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$wsBool.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$wsBool.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$wsBool.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$wsBool.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# A handful of blank, formatted rows were left below the data (rows 33-38).
$wsBool.Rows("33:38").Font.Name = $wsBool.Range("A32").Font.Name
$wsBool.Rows("33:38").Font.Size = $wsBool.Range("A32").Font.Size
$wsBool.Rows("33:38").Font.Bold = $wsBool.Range("A32").Font.Bold

# Restore the view: scrolled down a bit, with A32 selected.
$wsBool.Activate()
$wsBool.Range("A10").Select()
$excel.ActiveWindow.ScrollRow = 10
$wsBool.Range("A32").Select()

# --- Integer sheet: no data changes, just remember the last selection. ---
$wsInt = $wb.Worksheets.Item("Integer")
$wsInt.Range("A13").Select()

# --- About sheet: becomes the active tab again. ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
